# Appends 5 new task rows (138-142) to the "SB Squares Tasks" sheet,
# matching the pattern of the existing task data rows (A: Task #, B: Subject,
# C: Type, D: Assigned To, E: Status, F: Files Changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 138; TaskNum = 126; Subject = "QR code session redirect — redirect returning players to game view"; Type = "Feature"; Assignee = "ui-dev"; Status = "Completed"; Files = "app/join/[gameCode]/page.tsx" },
    @{ Row = 139; TaskNum = 127; Subject = "Game code reminder toast after joining via QR"; Type = "Feature"; Assignee = "ui-dev"; Status = "Completed"; Files = "app/join/[gameCode]/page.tsx, app/game/[gameId]/page.tsx" },
    @{ Row = 140; TaskNum = 128; Subject = "Redirect /join to /?code= when game not accepting players"; Type = "Feature"; Assignee = "ui-dev"; Status = "Completed"; Files = "app/join/[gameCode]/page.tsx, app/page.tsx" },
    @{ Row = 141; TaskNum = 129; Subject = "Add system theme option with 3-segment toggle (dark/light/system)"; Type = "Feature"; Assignee = "architect"; Status = "Completed"; Files = "hooks/use-theme.ts, app/layout.tsx, app/game/[gameId]/page.tsx" },
    @{ Row = 142; TaskNum = 130; Subject = "Improve grid cell border visibility in dark mode"; Type = "Enhancement"; Assignee = "architect"; Status = "Completed"; Files = "components/GridCell.tsx, components/Grid.tsx, lib/utils.ts" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.TaskNum
    $ws.Cells.Item($r.Row, 2).Value = $r.Subject
    $ws.Cells.Item($r.Row, 3).Value = $r.Type
    $ws.Cells.Item($r.Row, 4).Value = $r.Assignee
    $ws.Cells.Item($r.Row, 5).Value = $r.Status
    $ws.Cells.Item($r.Row, 6).Value = $r.Files
}

Write-Host "Added rows 138-142"
